# Insert a new data row at row 85 (pushing existing rows 85-139 down to 86-140)
# and populate it with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 85, shifting rows 85..139 down to 86..140.
$ws.Rows.Item(85).Insert()

# Populate the newly inserted row 85 with the new record's data.
$ws.Cells.Item(85, 1).Value = 3
$ws.Cells.Item(85, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(85, 3).Value = "Coquimbo"
$ws.Cells.Item(85, 4).Value = 44606
$ws.Cells.Item(85, 5).Value = 5
$ws.Cells.Item(85, 6).Value = 100112052
$ws.Cells.Item(85, 7).Value = "Albahaca"
$ws.Cells.Item(85, 8).Value = "Sin especificar"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 120
$ws.Cells.Item(85, 11).Value = 4500
$ws.Cells.Item(85, 12).Value = 5000
$ws.Cells.Item(85, 13).Value = 4750
$ws.Cells.Item(85, 14).Value = "`$/docena de matas"
$ws.Cells.Item(85, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(85, 16).Value = 792
$ws.Cells.Item(85, 17).Value = 6
$ws.Cells.Item(85, 18).Value = "Hortaliza"
